$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally had:
#   Row1: column headers (구분/직위/성명/주민등록번호/교육이수번호.../경력시작일/경력종료일/상근여부/전문인력여부)
#   Row2: a "placeholder / legend" sample row (신규=1 경력=2, 대리, 신동환, 000000-0000000, YYYY-MM-DD ...)
#   Row3: the actual sample data row (과장, 신동환, 880131-0000000, ...)
#
# The edit folds the legend text from row 2 into the row-1 header captions
# (구분\n신규 = 1\n경력 = 2, 주민등록번호\n000000-0000000, ...) and removes the
# now-redundant row 2, leaving a single data row (updated with new sample
# values) right under the header.

# 1) Remove the old legend/placeholder row - this shifts the old row 3
#    (the real sample data row) up into row 2.
$ws.Rows("2:2").Delete()

# 2) Fold the legend text into the header captions on row 1.
$ws.Range("A1").Value = "구분`r`n신규 = 1`r`n경력 = 2"
$ws.Range("D1").Value = "주민등록번호`r`n000000-0000000"
$ws.Range("E1").Value = "교육이수번호/인증서번호`r`n0000000000"
$ws.Range("F1").Value = "경력시작일`r`nYYYY-MM-DD"
$ws.Range("G1").Value = "경력종료일`r`nYYYY-MM-DD"
$ws.Range("H1").Value = "상근여부`r`n상근 = 1`r`n비상근 = 2"
$ws.Range("I1").Value = "전문인력여부`r`n전문인력 = 1`r`n비전문인력 = 2"

# 3) Wrap the now multi-line header captions and grow the header row to fit.
$ws.Range("A1").WrapText = $true
$ws.Range("D1:I1").WrapText = $true
$ws.Rows("1:1").RowHeight = 52.2

# 4) Refresh the sample data row with the new values.
$ws.Range("C2").Value = "홍길동"
$ws.Range("E2").Value = 2021051312
$ws.Range("F2").Value = "2020-04-01"
$ws.Range("G2").Value = "2023-08-15"

# 5) Move the active selection like the original author left it.
$ws.Range("E8").Select() | Out-Null
